$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Supplier and cost")

# New supplier rows to append to the "Supplier and cost" table
$ws.Range("A6").Value = "Brick"
$ws.Range("B6").Value = "Mr Brick"
$ws.Range("C6").Value = 20

$ws.Range("A7").Value = "Glass"
$ws.Range("B7").Value = "The GlassMaster"
$ws.Range("C7").Value = 50

$ws.Range("A8").Value = "Metal"
$ws.Range("B8").Value = "PedalToTheMetal"
$ws.Range("C8").Value = 100

$ws.Range("A9").Value = "Metal"
$ws.Range("B9").Value = "Metallica"
$ws.Range("C9").Value = 76

# Resize the supplier/cost table to include the newly added rows
$ws.ListObjects.Item("Table2").Resize($ws.Range("A1:C9"))

# Select the cell below the new data, matching the final selection state
$ws.Range("C10").Select()

# Make the "Supplier and cost" sheet the active tab
$ws.Activate()
